$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("August")

$ws.Range("B2").Value = 1923
$ws.Range("C2").Value = 1314
$ws.Range("D2").Value = 609
$ws.Range("E2").Value = "We borrowerd more than we lent"
$ws.Range("G2").Value = "1.46 : 1"

$ws.Range("B3").Value = 647
$ws.Range("C3").Value = 549
$ws.Range("D3").Value = 98
$ws.Range("E3").Value = "We borrowerd more than we lent"
$ws.Range("G3").Value = "1.18 : 1"

$ws.Range("B4").Value = 1397
$ws.Range("C4").Value = 1339
$ws.Range("D4").Value = 58
$ws.Range("E4").Value = "We borrowerd more than we lent"
$ws.Range("G4").Value = "1.04 : 1"

$ws.Range("B5").Value = 73
$ws.Range("C5").Value = 150
$ws.Range("D5").Value = -77
$ws.Range("F5").Value = "We lent more than we borrowed"
$ws.Range("G5").Value = "0.49 : 1"

$ws.Range("B6").Value = 1247
$ws.Range("C6").Value = 1654
$ws.Range("D6").Value = -407
$ws.Range("F6").Value = "We lent more than we borrowed"
$ws.Range("G6").Value = "0.75 : 1"

$ws.Range("B7").Value = 148
$ws.Range("C7").Value = 272
$ws.Range("D7").Value = -124
$ws.Range("F7").Value = "We lent more than we borrowed"
$ws.Range("G7").Value = "0.54 : 1"

$ws.Range("B8").Value = 148
$ws.Range("C8").Value = 187
$ws.Range("D8").Value = -39
$ws.Range("F8").Value = "We lent more than we borrowed"
$ws.Range("G8").Value = "0.79 : 1"

$ws.Range("B9").Value = 33
$ws.Range("C9").Value = 82
$ws.Range("D9").Value = -49
$ws.Range("F9").Value = "We lent more than we borrowed"
$ws.Range("G9").Value = "0.40 : 1"

$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 45
$ws.Range("D10").Value = -45
$ws.Range("F10").Value = "We lent more than we borrowed"
$ws.Range("G10").Value = "0.00 : 1"

$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 0

$ws.Range("B12").Value = 12
$ws.Range("C12").Value = 30
$ws.Range("D12").Value = -18
$ws.Range("F12").Value = "We lent more than we borrowed"
$ws.Range("G12").Value = "0.40 : 1"

$ws.Range("B13").Value = 201
$ws.Range("C13").Value = 98
$ws.Range("D13").Value = 103
$ws.Range("E13").Value = "We borrowerd more than we lent"
$ws.Range("G13").Value = "2.05 : 1"

$ws.Range("B14").Value = 142
$ws.Range("C14").Value = 315
$ws.Range("D14").Value = -173
$ws.Range("F14").Value = "We lent more than we borrowed"
$ws.Range("G14").Value = "0.45 : 1"

$ws.Range("B15").Value = 105
$ws.Range("C15").Value = 121
$ws.Range("D15").Value = -16
$ws.Range("F15").Value = "We lent more than we borrowed"
$ws.Range("G15").Value = "0.87 : 1"

$ws.Range("B16").Value = 55
$ws.Range("C16").Value = 178
$ws.Range("D16").Value = -123
$ws.Range("F16").Value = "We lent more than we borrowed"
$ws.Range("G16").Value = "0.31 : 1"

$ws.Range("B17").Value = 650
$ws.Range("C17").Value = 519
$ws.Range("D17").Value = 131
$ws.Range("E17").Value = "We borrowerd more than we lent"
$ws.Range("G17").Value = "1.25 : 1"

$ws.Range("B18").Value = 22
$ws.Range("C18").Value = 106
$ws.Range("D18").Value = -84
$ws.Range("F18").Value = "We lent more than we borrowed"
$ws.Range("G18").Value = "0.21 : 1"

$ws.Range("B19").Value = 638
$ws.Range("C19").Value = 547
$ws.Range("D19").Value = 91
$ws.Range("E19").Value = "We borrowerd more than we lent"
$ws.Range("G19").Value = "1.17 : 1"

$ws.Range("B20").Value = 3
$ws.Range("C20").Value = 90
$ws.Range("D20").Value = -87
$ws.Range("F20").Value = "We lent more than we borrowed"
$ws.Range("G20").Value = "0.03 : 1"

$ws.Range("B21").Value = 516
$ws.Range("C21").Value = 472
$ws.Range("D21").Value = 44
$ws.Range("E21").Value = "We borrowerd more than we lent"
$ws.Range("G21").Value = "1.09 : 1"

$ws.Range("B22").Value = 34
$ws.Range("C22").Value = 53
$ws.Range("D22").Value = -19
$ws.Range("F22").Value = "We lent more than we borrowed"
$ws.Range("G22").Value = "0.64 : 1"

$ws.Range("B23").Value = 726
$ws.Range("C23").Value = 436
$ws.Range("D23").Value = 290
$ws.Range("E23").Value = "We borrowerd more than we lent"
$ws.Range("G23").Value = "1.67 : 1"

$ws.Range("B24").Value = 1822
$ws.Range("C24").Value = 1430
$ws.Range("D24").Value = 392
$ws.Range("E24").Value = "We borrowerd more than we lent"
$ws.Range("G24").Value = "1.27 : 1"

$ws.Range("B25").Value = 128
$ws.Range("C25").Value = 378
$ws.Range("D25").Value = -250
$ws.Range("F25").Value = "We lent more than we borrowed"
$ws.Range("G25").Value = "0.34 : 1"

$ws.Range("B26").Value = 0
$ws.Range("C26").Value = 0
$ws.Range("D26").Value = 0

$ws.Range("B27").Value = 300
$ws.Range("C27").Value = 224
$ws.Range("D27").Value = 76
$ws.Range("E27").Value = "We borrowerd more than we lent"
$ws.Range("G27").Value = "1.34 : 1"

$ws.Range("B28").Value = 66
$ws.Range("C28").Value = 49
$ws.Range("D28").Value = 17
$ws.Range("E28").Value = "We borrowerd more than we lent"
$ws.Range("G28").Value = "1.35 : 1"

$ws.Range("B29").Value = 539
$ws.Range("C29").Value = 488
$ws.Range("D29").Value = 51
$ws.Range("E29").Value = "We borrowerd more than we lent"
$ws.Range("G29").Value = "1.10 : 1"

$ws.Range("B30").Value = 55
$ws.Range("C30").Value = 35
$ws.Range("D30").Value = 20
$ws.Range("E30").Value = "We borrowerd more than we lent"
$ws.Range("G30").Value = "1.57 : 1"

$ws.Range("B31").Value = 68
$ws.Range("C31").Value = 319
$ws.Range("D31").Value = -251
$ws.Range("F31").Value = "We lent more than we borrowed"
$ws.Range("G31").Value = "0.21 : 1"

$ws.Range("B32").Value = 539
$ws.Range("C32").Value = 540
$ws.Range("D32").Value = -1
$ws.Range("F32").Value = "We lent more than we borrowed"
$ws.Range("G32").Value = "1.00 : 1"

$ws.Range("B33").Value = 456
$ws.Range("C33").Value = 579
$ws.Range("D33").Value = -123
$ws.Range("F33").Value = "We lent more than we borrowed"
$ws.Range("G33").Value = "0.79 : 1"

$ws.Range("B34").Value = 226
$ws.Range("C34").Value = 119
$ws.Range("D34").Value = 107
$ws.Range("E34").Value = "We borrowerd more than we lent"
$ws.Range("G34").Value = "1.90 : 1"

$ws.Range("B35").Value = 1238
$ws.Range("C35").Value = 989
$ws.Range("D35").Value = 249
$ws.Range("E35").Value = "We borrowerd more than we lent"
$ws.Range("G35").Value = "1.25 : 1"

$ws.Range("B36").Value = 209
$ws.Range("C36").Value = 461
$ws.Range("D36").Value = -252
$ws.Range("F36").Value = "We lent more than we borrowed"
$ws.Range("G36").Value = "0.45 : 1"

$ws.Range("B37").Value = 502
$ws.Range("C37").Value = 362
$ws.Range("D37").Value = 140
$ws.Range("E37").Value = "We borrowerd more than we lent"
$ws.Range("G37").Value = "1.39 : 1"

$ws.Range("B38").Value = 13
$ws.Range("C38").Value = 177
$ws.Range("D38").Value = -164
$ws.Range("F38").Value = "We lent more than we borrowed"
$ws.Range("G38").Value = "0.07 : 1"

$ws.Range("B39").Value = 9
$ws.Range("C39").Value = 28
$ws.Range("D39").Value = -19
$ws.Range("F39").Value = "We lent more than we borrowed"
$ws.Range("G39").Value = "0.32 : 1"

$ws.Range("B40").Value = 24
$ws.Range("C40").Value = 49
$ws.Range("D40").Value = -25
$ws.Range("F40").Value = "We lent more than we borrowed"
$ws.Range("G40").Value = "0.49 : 1"

$ws.Range("B41").Value = 1
$ws.Range("C41").Value = 13
$ws.Range("D41").Value = -12
$ws.Range("F41").Value = "We lent more than we borrowed"
$ws.Range("G41").Value = "0.08 : 1"

$ws.Range("B42").Value = 19
$ws.Range("C42").Value = 11
$ws.Range("D42").Value = 8
$ws.Range("E42").Value = "We borrowerd more than we lent"
$ws.Range("G42").Value = "1.73 : 1"

$ws.Range("B43").Value = 0
$ws.Range("C43").Value = 0
$ws.Range("D43").Value = 0

$ws.Range("B44").Value = 70
$ws.Range("C44").Value = 55
$ws.Range("D44").Value = 15
$ws.Range("E44").Value = "We borrowerd more than we lent"
$ws.Range("G44").Value = "1.27 : 1"

$ws.Range("B45").Value = 70
$ws.Range("C45").Value = 125
$ws.Range("D45").Value = -55
$ws.Range("F45").Value = "We lent more than we borrowed"
$ws.Range("G45").Value = "0.56 : 1"

$ws.Range("B46").Value = 543
$ws.Range("C46").Value = 651
$ws.Range("D46").Value = -108
$ws.Range("F46").Value = "We lent more than we borrowed"
$ws.Range("G46").Value = "0.83 : 1"

$ws.Range("B47").Value = 859
$ws.Range("C47").Value = 591
$ws.Range("D47").Value = 268
$ws.Range("E47").Value = "We borrowerd more than we lent"
$ws.Range("G47").Value = "1.45 : 1"

$ws.Range("B48").Value = 227
$ws.Range("C48").Value = 654
$ws.Range("D48").Value = -427
$ws.Range("F48").Value = "We lent more than we borrowed"
$ws.Range("G48").Value = "0.35 : 1"

$ws.Range("B49").Value = 394
$ws.Range("C49").Value = 272
$ws.Range("D49").Value = 122
$ws.Range("E49").Value = "We borrowerd more than we lent"
$ws.Range("G49").Value = "1.45 : 1"

$ws.Range("B50").Value = 991
$ws.Range("C50").Value = 653
$ws.Range("D50").Value = 338
$ws.Range("E50").Value = "We borrowerd more than we lent"
$ws.Range("G50").Value = "1.52 : 1"

$ws.Range("B51").Value = 246
$ws.Range("C51").Value = 154
$ws.Range("D51").Value = 92
$ws.Range("E51").Value = "We borrowerd more than we lent"
$ws.Range("G51").Value = "1.60 : 1"

$ws.Range("B52").Value = 347
$ws.Range("C52").Value = 514
$ws.Range("D52").Value = -167
$ws.Range("F52").Value = "We lent more than we borrowed"
$ws.Range("G52").Value = "0.68 : 1"

$ws.Range("B53").Value = 154
$ws.Range("C53").Value = 257
$ws.Range("D53").Value = -103
$ws.Range("F53").Value = "We lent more than we borrowed"
$ws.Range("G53").Value = "0.60 : 1"

$ws.Range("B54").Value = 32
$ws.Range("C54").Value = 199
$ws.Range("D54").Value = -167
$ws.Range("F54").Value = "We lent more than we borrowed"
$ws.Range("G54").Value = "0.16 : 1"

$ws.Range("B55").Value = 259
$ws.Range("C55").Value = 193
$ws.Range("D55").Value = 66
$ws.Range("E55").Value = "We borrowerd more than we lent"
$ws.Range("G55").Value = "1.34 : 1"
